# Apply the "9 April 2020" New York hospitalization update:
#  - add a new "Deaths" column (F) with a header and historical back-fill for
#    the last several rows
#  - append a brand new row for 2020-04-09
#  - shrink the sheet-tab ratio in the workbook view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) New "Deaths" column header in F1, matching the style of the other
#        header cells (bold header font reuses the existing header style). ---
$ws.Range("F1").Value2 = "Deaths"
$ws.Range("F1").Font.Bold = $true

# --- 2) Back-fill Deaths values for the existing rows 21-25 (column F),
#        reusing the plain numeric style already used by column E (rows 24
#        and 25 have no E cell to copy from, so pull the shared plain-number
#        format from E21 instead - same underlying style index either way). ---
$ws.Range("E21").Copy($ws.Range("F21"))
$ws.Range("F21").Value2 = 594

$ws.Range("E22").Copy($ws.Range("F22"))
$ws.Range("F22").Value2 = 599

$ws.Range("E23").Copy($ws.Range("F23"))
$ws.Range("F23").Value2 = 731

$ws.Range("E21").Copy($ws.Range("F24"))
$ws.Range("F24").Value2 = 779

$ws.Range("E21").Copy($ws.Range("F25"))
$ws.Range("F25").Value2 = 799

# --- 3) Append the new row 26 for 2020-04-09, copying formats forward from
#        row 25 (date style for A, plain numeric for B:D). ---
$ws.Range("A25:D25").Copy($ws.Range("A26:D26"))
$ws.Range("A26").Value2 = 43930
$ws.Range("B26").Value2 = 290
$ws.Range("C26").Value2 = -17
$ws.Range("D26").Value2 = 109

$ws.Range("E21").Copy($ws.Range("F26"))
$ws.Range("F26").Value2 = 777

# --- 4) Shrink the tab-bar / horizontal-scrollbar split ratio (995 -> 500,
#        on Excel's internal 0..1650 TabRatio scale used by the COM model). ---
$win = $excel.ActiveWindow
$win.TabRatio = 500 / 1650
